$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format for the new row 9 to match existing data rows (style index 1, format "0")
$ws.Range("A9:W9").NumberFormat = $ws.Range("A8:W8").NumberFormat

# Recalculated values for existing rows (2-8), columns K:W

# Row 2
$ws.Range("K2").Value = 512834.50193147361
$ws.Range("L2").Value = 29525.700039975345
$ws.Range("M2").Value = 286425.10028348118
$ws.Range("N2").Value = 288156.40013650805
$ws.Range("O2").Value = 202032.30010768771
$ws.Range("P2").Value = 34558.700001463294
$ws.Range("Q2").Value = 282267.70011635125
$ws.Range("R2").Value = 277841.40023653209
$ws.Range("S2").Value = 24070.700017467141
$ws.Range("T2").Value = 315950.80032345653
$ws.Range("U2").Value = 490188.70024419576
$ws.Range("V2").Value = 316826.40011781454
$ws.Range("W2").Value = 301912.10025399923

# Row 3
$ws.Range("K3").Value = 512834.50193147361
$ws.Range("L3").Value = 29475.300040096045
$ws.Range("M3").Value = 286085.40028204024
$ws.Range("N3").Value = 288418.50013685226
$ws.Range("O3").Value = 202101.10010831058
$ws.Range("P3").Value = 34556.400001585484
$ws.Range("Q3").Value = 282350.90011708438
$ws.Range("R3").Value = 277817.00023600459
$ws.Range("S3").Value = 24073.400017492473
$ws.Range("T3").Value = 315560.70032213628
$ws.Range("U3").Value = 490519.60024516284
$ws.Range("V3").Value = 316907.30011866987
$ws.Range("W3").Value = 301890.40025349706

# Row 4
$ws.Range("K4").Value = 512834.50193147361
$ws.Range("L4").Value = 29390.400040067732
$ws.Range("M4").Value = 285562.30028200895
$ws.Range("N4").Value = 288600.80013689399
$ws.Range("O4").Value = 202185.00010865927
$ws.Range("P4").Value = 34562.800001643598
$ws.Range("Q4").Value = 282504.60012017936
$ws.Range("R4").Value = 277994.00023248047
$ws.Range("S4").Value = 24078.100017532706
$ws.Range("T4").Value = 314952.70032207668
$ws.Range("U4").Value = 490785.80024555326
$ws.Range("V4").Value = 317067.40012182295
$ws.Range("W4").Value = 302072.10025001317

# Row 5
$ws.Range("K5").Value = 512834.50193147361
$ws.Range("L5").Value = 29275.500041306019
$ws.Range("M5").Value = 285441.20028055459
$ws.Range("N5").Value = 288864.70013753325
$ws.Range("O5").Value = 202132.40010667592
$ws.Range("P5").Value = 34566.700001627207
$ws.Range("Q5").Value = 282391.90012191236
$ws.Range("R5").Value = 278127.50023224205
$ws.Range("S5").Value = 24078.100017614663
$ws.Range("T5").Value = 314716.70032186061
$ws.Range("U5").Value = 490997.10024420917
$ws.Range("V5").Value = 316958.60012353957
$ws.Range("W5").Value = 302205.60024985671

# Row 6
$ws.Range("K6").Value = 512834.50193147361
$ws.Range("L6").Value = 29096.800042673945
$ws.Range("M6").Value = 285029.1002824977
$ws.Range("N6").Value = 288846.10013855249
$ws.Range("O6").Value = 202297.70010128617
$ws.Range("P6").Value = 34570.200001642108
$ws.Range("Q6").Value = 282606.20012342185
$ws.Range("R6").Value = 278355.50023179501
$ws.Range("S6").Value = 24076.400017596781
$ws.Range("T6").Value = 314125.90032517165
$ws.Range("U6").Value = 491143.80023983866
$ws.Range("V6").Value = 317176.40012506396
$ws.Range("W6").Value = 302431.90024939179

# Row 7
$ws.Range("K7").Value = 512834.50193147361
$ws.Range("L7").Value = 28955.700042121112
$ws.Range("M7").Value = 284607.70027782023
$ws.Range("N7").Value = 289113.80014310777
$ws.Range("O7").Value = 202300.50010088086
$ws.Range("P7").Value = 34561.400001622736
$ws.Range("Q7").Value = 282817.50012075156
$ws.Range("R7").Value = 278447.40023555607
$ws.Range("S7").Value = 24074.000017605722
$ws.Range("T7").Value = 313563.40031994134
$ws.Range("U7").Value = 491414.30024398863
$ws.Range("V7").Value = 317378.9001223743
$ws.Range("W7").Value = 302521.40025316179

# Row 8
$ws.Range("K8").Value = 512834.50193147361
$ws.Range("L8").Value = 28868.400040626526
$ws.Range("M8").Value = 284572.70027856529
$ws.Range("N8").Value = 289140.70013933629
$ws.Range("O8").Value = 202377.70010089874
$ws.Range("P8").Value = 34567.100001655519
$ws.Range("Q8").Value = 282806.50012344867
$ws.Range("R8").Value = 278472.10023744404
$ws.Range("S8").Value = 24072.800017490983
$ws.Range("T8").Value = 313441.10031919181
$ws.Range("U8").Value = 491518.40024023503
$ws.Range("V8").Value = 317373.60012510419
$ws.Range("W8").Value = 302544.90025493503

# Row 9
$ws.Range("A9").Value = 2015
$ws.Range("B9").Value = 17914.399990998209
$ws.Range("C9").Value = 366334.10001897067
$ws.Range("D9").Value = 404582.60123407841
$ws.Range("E9").Value = 413558.00026299059
$ws.Range("F9").Value = 121249.89995732903
$ws.Range("G9").Value = 403148.49999554455
$ws.Range("H9").Value = 66902.600492224097
$ws.Range("I9").Value = 92037.600671075284
$ws.Range("J9").Value = 51984.800247728825
$ws.Range("K9").Value = 512834.50193147361
$ws.Range("L9").Value = 28808.400041893125
$ws.Range("M9").Value = 284590.40027844906
$ws.Range("N9").Value = 289152.90013868362
$ws.Range("O9").Value = 202421.70010026544
$ws.Range("P9").Value = 34567.400001659989
$ws.Range("Q9").Value = 282727.60012447834
$ws.Range("R9").Value = 278535.3002364859
$ws.Range("S9").Value = 24074.300017550588
$ws.Range("T9").Value = 313398.80032034218
$ws.Range("U9").Value = 491574.60023894906
$ws.Range("V9").Value = 317295.00012613833
$ws.Range("W9").Value = 302609.60025403649
